$wb = $excel.ActiveWorkbook

$wsTemplete = $wb.Worksheets.Item("templete")
$wsFyit = $wb.Worksheets.Item("fyit")

# templete sheet: A3 now holds the "sheetName" label (previously blank, kept its fill style)
$wsTemplete.Range("A3").Value = "sheetName"

# Move the active selection on templete to A3
$wsTemplete.Range("A3").Select()

# fyit sheet: update numeric values in column A
$wsFyit.Range("A2").Value = 2
$wsFyit.Range("A3").Value = 2.2000000000000002

# Make "fyit" the active sheet/tab (second sheet, index 1 zero-based)
$wsFyit.Activate()
$wsFyit.Range("A2").Select()

$wb.Save()
